# Add the 2020 data column (Q) to the SDG 1.3.1.1f indicator table.
# Formats are cloned from the existing 2019 column (P) so the new cells
# pick up the same fonts/number-formats/borders already used by the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (thin divider row under the title) - same blank formatted cell as P2.
$ws.Range("P2").Copy()
$ws.Range("Q2").PasteSpecial(-4122)

# Row 3 (year header) - 2020.
$ws.Range("P3").Copy()
$ws.Range("Q3").PasteSpecial(-4122)
$ws.Range("Q3").Value = 2020

# Row 4 (population receiving pensions/disability benefits).
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 197792

# Row 5 (share of total population, percent).
$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 2.9794303052841493

# Match the selection state stored in the saved file.
$ws.Range("G27").Select() | Out-Null
